$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 cells hold numeric-looking text (thousands separators etc.) that must
# stay as literal text, not be auto-converted to numbers by Excel's "smart"
# input parsing. Temporarily force a text number format so the assignment is
# stored verbatim, then clear the formatting again so no stray style is left
# on the cells (matches the original unstyled row 2 cells).
$rng = $ws.Range("A2:E2")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = ".00"
$ws.Range("B2").Value = "545,474,228.14"
$ws.Range("C2").Value = "1,570.00"
$ws.Range("D2").Value = "545,472,658.14"
$ws.Range("E2").Value = "22,416.14"

$rng.ClearFormats()
